$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.770.76"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -5.17%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.811.65"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -4.43%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "276.78"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -9.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5096"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -5.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3522"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -6.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.80"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06663"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -8.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.06"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -8.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.8353"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -6.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07837"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.793.32"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -4.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.072"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -4.86%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.79"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -7.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.000"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.20%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -6.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008007"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -7.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9994"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "25.848.53"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.735"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -5.77%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -6.96%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.060"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -6.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.207"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.42"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.656"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.02"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -6.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "108.92"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -6.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.347"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -9.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.221"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -9.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08793"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04893"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7343"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -10.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.138"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.894"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -4.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9993"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.051"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -7.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5235"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -11.84%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01854"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -6.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.287"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -13.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9520"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -11.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "111.94"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.63%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -6.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.137"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -11.78%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9997"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4579"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -9.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1381"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -9.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.299"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -8.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.30"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.501"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -7.47%  "
